# Desglose y creacion de Testcases variados
# Splits the single "CP003" sample-data row into its own row, duplicates it
# to create CP004 / CP005, and tweaks a couple of values on the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DatosCP")

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O")

# --- 1. Row 4 keeps the existing CP003 values (moved out of row 2). ---
$ws.Range("B4").Value = "Nombre"
$ws.Range("C4").Value = "Apellido"
$ws.Range("D4").Value = "DNI"
$ws.Range("E4").Value = 66565424
$ws.Range("F4").Value = "MASCULINO"
$ws.Range("G4").Value = "16-9-2005"
$ws.Range("H4").Value = "BC15Ejemplo@gmail.com"
$ws.Range("I4").Value = "BC15Ejemplo@gmail.com"
$ws.Range("J4").Value = "Ciudad de Buenos Aires"
$ws.Range("K4").Value = "Parque Chas"
$ws.Range("L4").Value = "Avenida de ejemplo 1234"
$ws.Range("M4").Value = 1123456789
$ws.Range("N4").Value = "Argentina"
$ws.Range("O4").Value = 14

# --- 2. Row 5: new test case CP004 (same sample data). ---
$ws.Range("A5").Value = "CP004"
$ws.Range("B5").Value = "Nombre"
$ws.Range("C5").Value = "Apellido"
$ws.Range("D5").Value = "DNI"
$ws.Range("E5").Value = 66565424
$ws.Range("F5").Value = "MASCULINO"
$ws.Range("G5").Value = "16-9-2005"
$ws.Range("H5").Value = "BC15Ejemplo@gmail.com"
$ws.Range("I5").Value = "BC15Ejemplo@gmail.com"
$ws.Range("J5").Value = "Ciudad de Buenos Aires"
$ws.Range("K5").Value = "Parque Chas"
$ws.Range("L5").Value = "Avenida de ejemplo 1234"
$ws.Range("M5").Value = 1123456789
$ws.Range("N5").Value = "Argentina"
$ws.Range("O5").Value = 14

# --- 3. Row 6: new test case CP005 (different DNI number / birth date). ---
$ws.Range("A6").Value = "CP005"
$ws.Range("B6").Value = "Nombre"
$ws.Range("C6").Value = "Apellido"
$ws.Range("D6").Value = "DNI"
$ws.Range("E6").Value = 66565426
$ws.Range("F6").Value = "MASCULINO"
$ws.Range("G6").Value = "16-9-2024"
$ws.Range("H6").Value = "BC15Ejemplo@gmail.com"
$ws.Range("I6").Value = "BC15Ejemplo@gmail.com"
$ws.Range("J6").Value = "Ciudad de Buenos Aires"
$ws.Range("K6").Value = "Parque Chas"
$ws.Range("L6").Value = "Avenida de ejemplo 1234"
$ws.Range("M6").Value = 1123456789
$ws.Range("N6").Value = "Argentina"
$ws.Range("O6").Value = 14

# --- 4. Clone the formatting of the sample-data row (currently row 2) onto
#        rows 4, 5 and 6 now that their values are in place. Copy/paste is
#        done per-cell, and *after* the value is written, so every column
#        keeps reusing its original style instead of Excel synthesizing a
#        new numFmt/style entry (this matters for the date column, which
#        carries a quoted-text + date numFmt combo). ---
foreach ($col in $cols) {
    $ws.Range($col + "2").Copy()
    $ws.Range($col + "4").PasteSpecial(-4122)
    $ws.Range($col + "2").Copy()
    $ws.Range($col + "5").PasteSpecial(-4122)
    $ws.Range($col + "2").Copy()
    $ws.Range($col + "6").PasteSpecial(-4122)
}

# --- 5. Row 2 only keeps the "CP001" label now; strip the sample data and
#        the tall row height that used to go with it. ---
$ws.Range("B2:O2").Clear()
$ws.Rows("2:2").AutoFit()

# --- 6. Drop the stray empty, styled D3 cell left over from the old layout. ---
$ws.Range("D3").Clear()

# --- 7. New data rows reuse the original ht="26" row height. ---
$ws.Rows("4:4").RowHeight = 26
$ws.Rows("5:5").RowHeight = 26
$ws.Rows("6:6").RowHeight = 26

# --- 8. Match the saved selection from the authored workbook. ---
$ws.Range("G7").Select()
